$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42; this shifts the existing rows 42:76 down to 43:77
# and expands the sheet dimension from A1:R76 to A1:R77 automatically.
$ws.Rows("42:42").Insert()

# Populate the newly inserted row 42 with the new record's data.
$ws.Cells.Item(42, 1).Value = 3
$ws.Cells.Item(42, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(42, 3).Value = "Coquimbo"
$ws.Cells.Item(42, 4).Value = 44907
$ws.Cells.Item(42, 5).Value = 5
$ws.Cells.Item(42, 6).Value = 100112022
$ws.Cells.Item(42, 7).Value = "Arveja Verde"
$ws.Cells.Item(42, 8).Value = "Perfection"
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 73
$ws.Cells.Item(42, 11).Value = 18000
$ws.Cells.Item(42, 12).Value = 19000
$ws.Cells.Item(42, 13).Value = 18521
$ws.Cells.Item(42, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(42, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(42, 16).Value = 741
$ws.Cells.Item(42, 17).Value = 25
$ws.Cells.Item(42, 18).Value = "Hortaliza"
